# Applies two changes to the document:
#  1. Split the run "left-top" into two runs: "left-" (keeping its
#     original character formatting) and "middle" (plain / no explicit
#     color), mirroring how Word splits a run when new text with
#     different direct formatting is inserted into the middle of it.
#  2. Explicitly set the footer distance on the section's page margins,
#     which serializes as w:footer="720" (720 twips = 36 pt) on
#     <w:pgMar>.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) left-top -> left- / middle
# ---------------------------------------------------------------------

# Build a brand-new run holding the replacement word "middle" in a
# location where it won't inherit any neighboring direct character
# formatting (color). We do this by appending a fresh paragraph at the
# very end of the story and typing into it - freshly typed text with no
# preceding run on an otherwise-empty paragraph picks up no explicit
# w:color, exactly matching the target run's <w:rPr><w:rtl w:val="0"/></w:rPr>.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$tempPara = $d.Paragraphs.Last
$tempPara.Range.InsertAfter("middle")
$tempPara = $d.Paragraphs.Last

# Range covering just the new word (excludes the paragraph mark).
$middleOnly = $d.Range($tempPara.Range.Start, $tempPara.Range.Start + 6)

# Stash its formatted content on the clipboard so it can be pasted
# in-place, preserving its (lack of) direct formatting.
$middleOnly.Copy()

# Locate "top" inside the original "left-top" run and overwrite it with
# the clipboard contents - this splits the original run into
# "left-" (unchanged formatting) followed by the newly pasted "middle"
# run (no color).
$target = $d.Content.Duplicate
$target.Find.Execute("top", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Paste()

# Remove the scratch paragraph we used to mint the clean "middle" run.
$scratch = $d.Paragraphs.Last
$scratch.Range.Delete()

# ---------------------------------------------------------------------
# 2) sectPr / pgMar footer distance
# ---------------------------------------------------------------------

# 36 points == 720 twips; setting this explicitly causes w:footer="720"
# to be written on <w:pgMar>.
$d.PageSetup.FooterDistance = 36
